# Apply data updates described by the diff across the four affected sheets:
# model_summary, score_psi_test_segments, feature_importance, woe_bins_feature_08
$wb = $excel.ActiveWorkbook

# --- model_summary (sheet1) ---
$ws = $wb.Worksheets("model_summary")
$ws.Range("A2").Value = "XGBoost"
$ws.Range("B2").Value = 0.8568094525436677

# --- score_psi_test_segments (sheet10) ---
$ws = $wb.Worksheets("score_psi_test_segments")
$ws.Range("B2").Value = "(-0.001, 0.0102]"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 11
$ws.Range("E2").Value = 400
$ws.Range("F2").Value = 110
$ws.Range("G2").Value = 0.9999999999999996
$ws.Range("H2").Value = 0.000953101798043247
$ws.Range("I2").Value = 0.005410111043602228
$ws.Range("J2").Value = 0.005115886218845844
$ws.Range("K2").Value = 0.000953101798043247
$ws.Range("B3").Value = "(0.0102, 0.0239]"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 9.700000000000001
$ws.Range("E3").Value = 400
$ws.Range("F3").Value = 97
$ws.Range("G3").Value = -0.3000000000000003
$ws.Range("H3").Value = [double]"9.13776224541258e-05"
$ws.Range("I3").Value = 0.01648161187767982
$ws.Range("J3").Value = 0.01650168187916279
$ws.Range("K3").Value = 0.001044479420497373
$ws.Range("B4").Value = "(0.0239, 0.0431]"
$ws.Range("D4").Value = 8.799999999999999
$ws.Range("F4").Value = 88
$ws.Range("G4").Value = -1.200000000000001
$ws.Range("H4").Value = 0.001534000458118622
$ws.Range("I4").Value = 0.03172575309872627
$ws.Range("J4").Value = 0.03246678039431572
$ws.Range("K4").Value = 0.002578479878615994
$ws.Range("B5").Value = "(0.0431, 0.0698]"
$ws.Range("C5").Value = 10.05
$ws.Range("D5").Value = 11.2
$ws.Range("E5").Value = 402
$ws.Range("F5").Value = 112
$ws.Range("G5").Value = 1.15
$ws.Range("H5").Value = 0.001245923153653587
$ws.Range("I5").Value = 0.05581361427903175
$ws.Range("J5").Value = 0.05511580035090446
$ws.Range("K5").Value = 0.003824403032269581
$ws.Range("B6").Value = "(0.0698, 0.121]"
$ws.Range("C6").Value = 9.950000000000001
$ws.Range("D6").Value = 9.800000000000001
$ws.Range("E6").Value = 398
$ws.Range("F6").Value = 98
$ws.Range("G6").Value = -0.1500000000000001
$ws.Range("H6").Value = [double]"2.278524824096287e-05"
$ws.Range("I6").Value = 0.09360076487064362
$ws.Range("J6").Value = 0.09032906591892242
$ws.Range("K6").Value = 0.003847188280510544
$ws.Range("B7").Value = "(0.121, 0.196]"
$ws.Range("C7").Value = 10.15
$ws.Range("D7").Value = 10.3
$ws.Range("E7").Value = 406
$ws.Range("F7").Value = 103
$ws.Range("G7").Value = 0.1499999999999987
$ws.Range("H7").Value = [double]"2.200528462169025e-05"
$ws.Range("I7").Value = 0.1576640456914902
$ws.Range("J7").Value = 0.1617581695318222
$ws.Range("K7").Value = 0.003869193565132234
$ws.Range("B8").Value = "(0.196, 0.306]"
$ws.Range("C8").Value = 9.875
$ws.Range("D8").Value = 9.4
$ws.Range("E8").Value = 395
$ws.Range("F8").Value = 94
$ws.Range("G8").Value = -0.4750000000000004
$ws.Range("H8").Value = 0.0002341589521783303
$ws.Range("I8").Value = 0.2552129626274109
$ws.Range("J8").Value = 0.253635436296463
$ws.Range("K8").Value = 0.004103352517310564
$ws.Range("B9").Value = "(0.306, 0.44]"
$ws.Range("C9").Value = 10.15
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = 406
$ws.Range("F9").Value = 80
$ws.Range("G9").Value = -2.15
$ws.Range("H9").Value = 0.005117691521871152
$ws.Range("I9").Value = 0.3768987655639648
$ws.Range("J9").Value = 0.3747358918190002
$ws.Range("K9").Value = 0.009221044039181716
$ws.Range("B10").Value = "(0.44, 0.686]"
$ws.Range("C10").Value = 9.825000000000001
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 393
$ws.Range("F10").Value = 110
$ws.Range("G10").Value = 1.175
$ws.Range("H10").Value = 0.001327340101755786
$ws.Range("I10").Value = 0.5581108331680298
$ws.Range("J10").Value = 0.5460732579231262
$ws.Range("K10").Value = 0.0105483841409375
$ws.Range("B11").Value = "(0.686, 1.0]"
$ws.Range("D11").Value = 10.8
$ws.Range("F11").Value = 108
$ws.Range("G11").Value = 0.7999999999999994
$ws.Range("H11").Value = 0.000615688329089025
$ws.Range("I11").Value = 0.822874903678894
$ws.Range("J11").Value = 0.8338835835456848
$ws.Range("K11").Value = 0.01116407247002653

# --- feature_importance (sheet2) ---
$ws = $wb.Worksheets("feature_importance")
$ws.Range("A2").Value = "feature_26"
$ws.Range("B2").Value = 0.1719878911972046
$ws.Range("A3").Value = "feature_46"
$ws.Range("B3").Value = 0.1675053536891937
$ws.Range("B4").Value = 0.1633881777524948
$ws.Range("B5").Value = 0.1343928128480911
$ws.Range("B6").Value = 0.1252419352531433
$ws.Range("B7").Value = 0.1250054389238358
$ws.Range("B8").Value = 0.1124783754348755

# --- woe_bins_feature_08 (sheet5) ---
$ws = $wb.Worksheets("woe_bins_feature_08")
$ws.Range("A2").Value = "[-12.02, -5.16]"
$ws.Range("B2").Value = -0.5098482432106972
$ws.Range("D2").Value = 374
$ws.Range("E2").Value = 0.06733167082294264
$ws.Range("F2").Value = 0.02138621439839053
$ws.Range("A3").Value = "[-5.16, -3.88]"
$ws.Range("B3").Value = -0.9829702172630708
$ws.Range("D3").Value = 382
$ws.Range("E3").Value = 0.04260651629072681
$ws.Range("F3").Value = 0.06598229373534219
$ws.Range("B4").Value = -0.2614816422153408
$ws.Range("C4").Value = 34
$ws.Range("D4").Value = 366
$ws.Range("E4").Value = 0.08500000000000001
$ws.Range("F4").Value = 0.006183460215060842
$ws.Range("A5").Value = "[-2.93, -2.12]"
$ws.Range("B5").Value = -0.6306873746152155
$ws.Range("D5").Value = 376
$ws.Range("E5").Value = 0.06
$ws.Range("F5").Value = 0.03114827565808744
$ws.Range("A6").Value = "[-2.12, -1.27]"
$ws.Range("B6").Value = -0.3606056394826176
$ws.Range("C6").Value = 31
$ws.Range("D6").Value = 369
$ws.Range("E6").Value = 0.0775
$ws.Range("F6").Value = 0.01131212396800935
$ws.Range("A7").Value = "[-1.27, -0.53]"
$ws.Range("B7").Value = -0.1096337297438762
$ws.Range("C7").Value = 39
$ws.Range("D7").Value = 360
$ws.Range("E7").Value = 0.09774436090225563
$ws.Range("F7").Value = 0.001150861669558705
$ws.Range("B8").Value = 0.0674846184240665
$ws.Range("C8").Value = 46
$ws.Range("D8").Value = 355
$ws.Range("E8").Value = 0.114713216957606
$ws.Range("F8").Value = 0.0004696824854299417
$ws.Range("B9").Value = 0.2304745439711655
$ws.Range("C9").Value = 53
$ws.Range("D9").Value = 347
$ws.Range("E9").Value = 0.1325
$ws.Range("F9").Value = 0.005821395375469932
$ws.Range("A10").Value = "[1.17, 2.65]"
$ws.Range("B10").Value = 0.2729479326187155
$ws.Range("C10").Value = 55
$ws.Range("D10").Value = 345
$ws.Range("E10").Value = 0.1375
$ws.Range("F10").Value = 0.008299339271090993
$ws.Range("A11").Value = "[2.65, 12.45]"
$ws.Range("B11").Value = 1.12244209854006
$ws.Range("C11").Value = 109
$ws.Range("D11").Value = 291
$ws.Range("E11").Value = 0.2725
$ws.Range("F11").Value = 0.1901445399464467

